# Append a new row (row 3) to the "ランサーズ" sheet with the 2026-02-07
# 12:43 JST scrape results, update the A2 timestamp to match the same
# run, widen column D slightly, and wire up the new row's hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: refresh the "取得日時" (fetched-at) timestamp for this run.
$ws.Range("A2").Value = "2026-02-07 12:43:02"

# Column D ("価格") needs a little more room for the new row's price text.
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668

# Row 3: newly scraped listing.
$ws.Range("A3").Value = "2026-02-07 12:43:02"
$ws.Range("B3").Value = "出品代行サービス用Webアプリ開発依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5487615"
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = "◆開発 ◇アプリ"

# F3 is a hyperlink (matching the style already used by F2).
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5487615")
$ws.Range("F3").Style = "Hyperlink"
